$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate row for "Carbon dioxide, non-fossil" / "air::non-urban air or from high stacks"
# which appears both as row 5 and row 8. Delete row 8 (the duplicate), shifting rows below up.
$ws.Rows.Item(8).Delete()
